$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The revision rewrites the "Objetivos/Programa resumido/Programa/Bibliografia" section
# (rows 10-21 originally) by inserting body paragraphs for several fields and shifting the
# remaining labels down, ending with one additional row (22) for "Bibliografia".
# The cleanest way to reproduce this with the object model is to clear the whole block
# and rewrite every cell with its final target content, copying the correct column
# formatting (bold header in col A, wrapped body text in col B/C with red font in C).

# --- Text content (defined first as here-strings for readability) ---

$t10A = @'
Objetivos:
'@

$t10B = @'
Proporcionar aos estudantes dos cursos de Engenharia da Escola de Engenharia de Lorena o contato com os fundamentos básicos das Ciências Humanas ou Sociais e estimular a reflexão sobre o desenvolvimento científico-tecnológico no mundo contemporâneo e seus reflexos na vida do homem no planeta e na sociedade brasileira em tempos de globalização.
Procurar assim, complementar a formação dos mesmos, capacitando-os para utilizar os conhecimentos adquiridos no seu contexto social, no exercício profissional, com competência, criatividade e com amplo entendimento da sua ação como cidadão responsável e solidário.
'@

$t10C = @'
Proporcionar aos estudantes dos cursos de Engenharia da Escola de Engenharia de Lorena o contato com os fundamentos básicos das Ciências Humanas ou Sociais e estimular a reflexão sobre o desenvolvimento científico-tecnológico no mundo contemporâneo e seus reflexos na vida do homem no planeta e na sociedade brasileira em tempos de globalização.
Procurar assim, complementar a formação dos mesmos, capacitando-os para utilizar os conhecimentos adquiridos no seu contexto social, no exercício profissional, com competência, criatividade e com amplo entendimento da sua ação como cidadão responsável e solidário.
'@

$t11A = @'
Objectives:
'@

$t12A = @'
Docentes responsáveis:
'@

$t13B = @'
6376612 - Daisy Rafaela da Silva
'@

$t13C = @'
6376612 - Daisy Rafaela da Silva
'@

$t14A = @'
Programa resumido:
'@

$t14B = @'
Introdução às Ciências Sociais. A sociedade do conhecimento. o homem na sociedade Emergente. Globalização e a realidade brasileira.
'@

$t14C = @'
Introdução às Ciências Sociais. A sociedade do conhecimento. o homem na sociedade Emergente. Globalização e a realidade brasileira.
'@

$t15A = @'
Short syllabus:
'@

$t16A = @'
Programa:
'@

$t16B = @'
Introdução às Ciências Sociais
- a posição das Ciências Sociais no quadro das ciências 
- o papel das Ciências Sociais na atualidade;
2 - A Sociedade do Conhecimento 
- A evolução do conhecimento 
- O paradigma científico e a revolução científica- tecnológica
- Mudanças no paradigma científico
- A questão sócio-ambiental
3 -  O homem na sociedade emergente
- ética nas relações humanas
- liderança pessoal e profissional
4 - Globalização e a realidade brasileira
- o sistema hegemônico : o neoliberalismo;
- a globalização econômica
- a globalização social 
- a sociedade civil globalizada
'@

$t16C = @'
Introdução às Ciências Sociais
- a posição das Ciências Sociais no quadro das ciências 
- o papel das Ciências Sociais na atualidade;
2 - A Sociedade do Conhecimento 
- A evolução do conhecimento 
- O paradigma científico e a revolução científica- tecnológica
- Mudanças no paradigma científico
- A questão sócio-ambiental
3 -  O homem na sociedade emergente
- ética nas relações humanas
- liderança pessoal e profissional
4 - Globalização e a realidade brasileira
- o sistema hegemônico : o neoliberalismo;
- a globalização econômica
- a globalização social 
- a sociedade civil globalizada
'@

$t17A = @'
Syllabus:
'@

$t18A = @'
Avaliação:
'@

$t19A = @'
Método:
'@

$t19B = @'
A média semestral e final dos alunos será composta por: Prova Semestral  (PS) e outros instrumentos (T) empregados na avaliação do aluno, valorizando a sua participação e colaboração nos trabalhos e atividades desenvolvidas individualmente e no Projeto de curso em equipe.
'@

$t19C = @'
A média semestral e final dos alunos será composta por: Prova Semestral  (PS) e outros instrumentos (T) empregados na avaliação do aluno, valorizando a sua participação e colaboração nos trabalhos e atividades desenvolvidas individualmente e no Projeto de curso em equipe.
'@

$t20A = @'
Critério:
'@

$t20B = @'
(PS+T) / 2
'@

$t20C = @'
(PS+T) / 2
'@

$t21A = @'
Norma de recuperação:
'@

$t21B = @'
-  Trabalho escrito, com questionamento, envolvendo o conteúdo do programa 
-   prova escrita
'@

$t21C = @'
-  Trabalho escrito, com questionamento, envolvendo o conteúdo do programa 
-   prova escrita
'@

$t22A = @'
Bibliografia:
'@

$t22B = @'
01  Constituição da República Federativa do Brasil
02  Código de Ética do Engenheiro, CREA: 2002.
03 -  CAPRA, F.   A Teia da Vida.  São Paulo: Cultrix, 2003.
03  CASTELLS, Manuel. O Poder da Identidade. A Era da Informação: Economia, Sociedade e Cultura. Vol. 2 . São Paulo: Paz e Terra, 1999.
04  - HUNTER, James C.  O Monge e o Executivo: uma história sobre a essência da liderança.  Rio de Janeiro: Sextante, 2004.
05  NOVAES, Adauto ( org.) Ética. São Paulo: Secretaria Municipal de Cultura e Companhia das Letras, 1992.
07  SADER, Emir. A Vingança da História. São Paulo: Boitempo-Editorial, 2003.
08  SANTOS, Boaventura de. Um Discurso sobre as Ciências. Porto, Portugal: Afrontamentos, 1997.
10 -  SCHAEFER, Richard T.  Sociologia. 6ª. Ed.; São Paulo: McGraw-Hill, 2006.
11  SODERO TOLEDO, Francisco. Outros Caminhos : Vale do Paraíba, do regional ao  internacional, do global ao local. São Paulo, Editora Salesiana, 2001.
_____________  Eu,Tu,Nós  Ética e Cidadania para jovens. Cachoeira Paulista, São Paulo: Ed. Canção Nova, 2005
Artigos de revistas especializadas e de jornais; 
Estudos, artigos, notícias e pesquisas via internet.
'@

$t22C = @'
01  Constituição da República Federativa do Brasil
02  Código de Ética do Engenheiro, CREA: 2002.
03 -  CAPRA, F.   A Teia da Vida.  São Paulo: Cultrix, 2003.
03  CASTELLS, Manuel. O Poder da Identidade. A Era da Informação: Economia, Sociedade e Cultura. Vol. 2 . São Paulo: Paz e Terra, 1999.
04  - HUNTER, James C.  O Monge e o Executivo: uma história sobre a essência da liderança.  Rio de Janeiro: Sextante, 2004.
05  NOVAES, Adauto ( org.) Ética. São Paulo: Secretaria Municipal de Cultura e Companhia das Letras, 1992.
07  SADER, Emir. A Vingança da História. São Paulo: Boitempo-Editorial, 2003.
08  SANTOS, Boaventura de. Um Discurso sobre as Ciências. Porto, Portugal: Afrontamentos, 1997.
10 -  SCHAEFER, Richard T.  Sociologia. 6ª. Ed.; São Paulo: McGraw-Hill, 2006.
11  SODERO TOLEDO, Francisco. Outros Caminhos : Vale do Paraíba, do regional ao  internacional, do global ao local. São Paulo, Editora Salesiana, 2001.
_____________  Eu,Tu,Nós  Ética e Cidadania para jovens. Cachoeira Paulista, São Paulo: Ed. Canção Nova, 2005
Artigos de revistas especializadas e de jornais; 
Estudos, artigos, notícias e pesquisas via internet.
'@

# --- Clear the existing block (values, number formats and styles) ---
$ws.Range("A10:C22").Clear() | Out-Null

# --- Template cells holding the correct formatting for each column ---
$headerTemplate = $ws.Cells.Item(3, 1)   # bold label style (col A)
$bodyTemplate   = $ws.Cells.Item(3, 2)   # wrapped body style (col B)
$redTemplate    = $ws.Cells.Item(3, 3)   # wrapped red body style (col C)

function Set-Cell($row, $col, $template, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $template.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
    $cell.Value = $value
}

# --- Row 10: Objetivos: / long objectives paragraph (same text in B and C) ---
Set-Cell 10 1 $headerTemplate $t10A
Set-Cell 10 2 $bodyTemplate $t10B
Set-Cell 10 3 $redTemplate $t10C

# --- Row 11: Objectives: ---
Set-Cell 11 1 $headerTemplate $t11A

# --- Row 12: Docentes responsaveis: ---
Set-Cell 12 1 $headerTemplate $t12A

# --- Row 13: (no label) professor name moved here ---
Set-Cell 13 2 $bodyTemplate $t13B
Set-Cell 13 3 $redTemplate $t13C

# --- Row 14: Programa resumido: / short syllabus paragraph ---
Set-Cell 14 1 $headerTemplate $t14A
Set-Cell 14 2 $bodyTemplate $t14B
Set-Cell 14 3 $redTemplate $t14C

# --- Row 15: Short syllabus: ---
Set-Cell 15 1 $headerTemplate $t15A

# --- Row 16: Programa: / full program text ---
Set-Cell 16 1 $headerTemplate $t16A
Set-Cell 16 2 $bodyTemplate $t16B
Set-Cell 16 3 $redTemplate $t16C

# --- Row 17: Syllabus: ---
Set-Cell 17 1 $headerTemplate $t17A

# --- Row 18: Avaliacao: ---
Set-Cell 18 1 $headerTemplate $t18A

# --- Row 19: Metodo: / evaluation method paragraph ---
Set-Cell 19 1 $headerTemplate $t19A
Set-Cell 19 2 $bodyTemplate $t19B
Set-Cell 19 3 $redTemplate $t19C

# --- Row 20: Criterio: / (PS+T)/2 ---
Set-Cell 20 1 $headerTemplate $t20A
Set-Cell 20 2 $bodyTemplate $t20B
Set-Cell 20 3 $redTemplate $t20C

# --- Row 21: Norma de recuperacao: / recovery rule paragraph ---
Set-Cell 21 1 $headerTemplate $t21A
Set-Cell 21 2 $bodyTemplate $t21B
Set-Cell 21 3 $redTemplate $t21C

# --- Row 22 (new row): Bibliografia: / bibliography list ---
Set-Cell 22 1 $headerTemplate $t22A
Set-Cell 22 2 $bodyTemplate $t22B
Set-Cell 22 3 $redTemplate $t22C

# --- Row heights: 60pt for short paragraphs, 120pt for the long ones, default otherwise ---
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).EntireRow.AutoFit() | Out-Null
$ws.Rows.Item(13).EntireRow.AutoFit() | Out-Null
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).EntireRow.AutoFit() | Out-Null
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120

